$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new rows of data (rows 11-13)
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = "8_Gate Set Tomography"
$ws.Range("E11").Value = "C:\Jeonghyun\GIT\QSCOUT\Gate_Set_Tomography"

$ws.Range("C12").Value = 9
$ws.Range("D12").Value = "9_Determination of Multi-mode Motional Quantum States in a Trapped Ion System"
$ws.Range("E12").Value = "C:\Jeonghyun\GIT\QSCOUT"

$ws.Range("C13").Value = 10
$ws.Range("D13").Value = "10_Characterizing and mitigating coherent errors in a trapped ion quantum processor using hidden inverses"
$ws.Range("E13").Value = "C:\Jeonghyun\GIT\QSCOUT"

# Adjust column widths to fit the new (longer) content
$ws.Columns.Item(4).ColumnWidth = 91.71428571428572
$ws.Columns.Item(5).ColumnWidth = 49

# Update the active selection to the last entered cell
$ws.Range("E13").Select()
